$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B4").Value = 8.639400000000002
$ws.Range("A9").Value = -22.10839999999999
$ws.Range("B9").Value = 6.602500000000005
$ws.Range("C9").Value = -12.0559
$ws.Range("B11").Value = 4.791700000000002
$ws.Range("A13").Value = -22.3101
$ws.Range("A16").Value = -21.69130000000001
$ws.Range("B16").Value = 4.808
$ws.Range("A18").Value = -22.15310000000001
$ws.Range("A20").Value = -21.42639999999998
$ws.Range("C22").Value = -12.2334
$ws.Range("B23").Value = 8.811199999999996
$ws.Range("B24").Value = 6.3044
$ws.Range("A26").Value = -21.10959999999997
$ws.Range("B26").Value = 5.530300000000002
$ws.Range("A27").Value = -21.93669999999999
$ws.Range("C27").Value = -12.66599999999999
$ws.Range("A29").Value = -21.81769999999998
$ws.Range("C29").Value = -11.88370000000001
$ws.Range("C32").Value = -12.6454
$ws.Range("B34").Value = 9.880000000000004
$ws.Range("A35").Value = -21.86579999999999
$ws.Range("B35").Value = 5.170800000000001
$ws.Range("A36").Value = -20.93599999999998
$ws.Range("C37").Value = -14.09099999999999
$ws.Range("C38").Value = -11.9268
$ws.Range("C39").Value = -12.70230000000001
$ws.Range("C41").Value = -12.97550000000002
$ws.Range("B44").Value = 4.512800000000002
$ws.Range("A45").Value = -21.62599999999999
$ws.Range("C45").Value = -13.37899999999999
$ws.Range("B48").Value = 7.002500000000001
$ws.Range("C48").Value = -12.24570000000001
$ws.Range("B49").Value = 5.987099999999998
$ws.Range("C51").Value = -11.12799999999999
$ws.Range("B52").Value = 5.5893
$ws.Range("A55").Value = -22.2021
$ws.Range("C56").Value = -12.74629999999999
$ws.Range("A57").Value = -21.97719999999999
$ws.Range("C57").Value = -12.56409999999998
$ws.Range("C61").Value = -14.18599999999999
$ws.Range("C64").Value = -10.10509999999999
$ws.Range("B66").Value = 4.966499999999995
$ws.Range("B67").Value = 5.142300000000001
$ws.Range("A69").Value = -21.5747
$ws.Range("B73").Value = 9.133800000000003
$ws.Range("C75").Value = -12.3128
$ws.Range("A76").Value = -19.96039999999999
$ws.Range("A78").Value = -21.63799999999999
$ws.Range("B78").Value = 5.858800000000002
$ws.Range("B80").Value = 9.314699999999998
$ws.Range("A82").Value = -21.9122
$ws.Range("C82").Value = -11.4159
$ws.Range("A83").Value = -21.60349999999999
$ws.Range("C90").Value = -10.0444
$ws.Range("B91").Value = 5.133799999999999
$ws.Range("A93").Value = -21.43600000000001
$ws.Range("C93").Value = -10.44889999999999
$ws.Range("A97").Value = -21.5742
$ws.Range("B97").Value = 5.065099999999995
$ws.Range("B99").Value = 5.316199999999998
$ws.Range("C102").Value = -11.6694
$ws.Range("B104").Value = 10.23450000000001
$ws.Range("C105").Value = -12.65190000000001
